# Apply "traded, fixed issues with the repeater":
#  - fill in PriceChange (X) / UpDown (Y) for the previous last row (row 10)
#  - append a new trade row (row 11) with the latest data point

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- back-fill the previous last row (10) now that the next day's data exists ---
$ws.Range("X10").Value = -0.010002000000000066
$ws.Range("Y10").Value = "Down"

# --- append the new row (11) for the latest trade/scrape ---
$ws.Range("A11").Value2 = 42654.882106481484
$ws.Range("B11").Value = 3
$ws.Range("C11").Value = "Neutral"
$ws.Range("D11").Value = 22
$ws.Range("E11").Value = 20034
$ws.Range("F11").Value = 3419
$ws.Range("G11").Value = 63
$ws.Range("H11").Value = 35
$ws.Range("I11").Value = 81
$ws.Range("J11").Value = 18
$ws.Range("K11").Value = 21576
$ws.Range("L11").Value = 394
$ws.Range("M11").Value = 224
$ws.Range("N11").Value = 104
$ws.Range("O11").Value = 24
$ws.Range("P11").Value = "Noun"
$ws.Range("Q11").Value = 17.089518681678967
$ws.Range("R11").Value = -24.44
$ws.Range("S11").Value = -0.1101
$ws.Range("T11").Value = -0.0419
$ws.Range("U11").Value = 6.47
$ws.Range("V11").Value = 1.88
$ws.Range("W11").Value = -2

# match formatting used by the rest of the "trade" rows
$ws.Range("A11").NumberFormat = "m/d/yy h:mm"
$ws.Range("S11").NumberFormat = "0.00%"
$ws.Range("T11").NumberFormat = "0.00%"
